$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "problem_user"
Write-Host "Sheets count: $($wb.Worksheets.Count)"
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    Write-Host "Sheet $i : $($wb.Worksheets.Item($i).Name)"
}
